$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on percentage cells to avoid Excel auto-converting
# "NN%" strings into numeric percentage values (which would change cell type/style).
$percentCells = @("H5", "H7", "H8", "H9", "H14", "H16", "H17", "H23", "H26", "H28", "H32", "H38", "H39", "H40", "H41", "H44", "H46")
foreach ($addr in $percentCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-11 18:48:33"
$ws.Range("I2").Value = "4.1 mm"
$ws.Range("E3").Value = "2026-02-11 18:48:35"
$ws.Range("I3").Value = "1.9 mm"
$ws.Range("O3").Value = "0.3 °C"
$ws.Range("E4").Value = "2026-02-11 18:48:38"
$ws.Range("J4").Value = "1002.5 hPa"
$ws.Range("O4").Value = "16.0 °C"
$ws.Range("E5").Value = "2026-02-11 18:48:41"
$ws.Range("H5").Value = "78%"
$ws.Range("I5").Value = "2.4 mm"
$ws.Range("E6").Value = "2026-02-11 18:48:43"
$ws.Range("J6").Value = "1003.1 hPa"
$ws.Range("E7").Value = "2026-02-11 18:48:46"
$ws.Range("H7").Value = "42%"
$ws.Range("N7").Value = "16.4 °C 18:29 TU"
$ws.Range("O7").Value = "19.3 °C"
$ws.Range("E8").Value = "2026-02-11 18:48:48"
$ws.Range("H8").Value = "53%"
$ws.Range("I8").Value = "0.4 mm"
$ws.Range("J8").Value = "1003.3 hPa"
$ws.Range("N8").Value = "12.0 °C 18:29 TU"
$ws.Range("O8").Value = "15.4 °C"
$ws.Range("E9").Value = "2026-02-11 18:48:51"
$ws.Range("H9").Value = "87%"
$ws.Range("E10").Value = "2026-02-11 18:48:53"
$ws.Range("E11").Value = "2026-02-11 18:48:56"
$ws.Range("I11").Value = "0.3 mm"
$ws.Range("E12").Value = "2026-02-11 18:48:58"
$ws.Range("E13").Value = "2026-02-11 18:49:00"
$ws.Range("J13").Value = "1005.3 hPa"
$ws.Range("E14").Value = "2026-02-11 18:49:03"
$ws.Range("H14").Value = "46%"
$ws.Range("K14").Value = "10.3 MJ/m2"
$ws.Range("E15").Value = "2026-02-11 18:49:05"
$ws.Range("E16").Value = "2026-02-11 18:49:08"
$ws.Range("H16").Value = "64%"
$ws.Range("I16").Value = "6.3 mm"
$ws.Range("E17").Value = "2026-02-11 18:49:10"
$ws.Range("H17").Value = "76%"
$ws.Range("E18").Value = "2026-02-11 18:49:13"
$ws.Range("J18").Value = "1003.1 hPa"
$ws.Range("L18").Value = "26.6 km/h - 273º 18:02 TU"
$ws.Range("O18").Value = "14.1 °C"
$ws.Range("E19").Value = "2026-02-11 18:49:15"
$ws.Range("E20").Value = "2026-02-11 18:49:18"
$ws.Range("I20").Value = "1.0 mm"
$ws.Range("E21").Value = "2026-02-11 18:49:20"
$ws.Range("I21").Value = "1.8 mm"
$ws.Range("J21").Value = "1005.8 hPa"
$ws.Range("E22").Value = "2026-02-11 18:49:23"
$ws.Range("O22").Value = "-2.7 °C"
$ws.Range("E23").Value = "2026-02-11 18:49:25"
$ws.Range("H23").Value = "72%"
$ws.Range("I23").Value = "3.6 mm"
$ws.Range("L23").Value = "71.3 km/h - 268º 18:03 TU"
$ws.Range("E24").Value = "2026-02-11 18:49:28"
$ws.Range("I24").Value = "7.5 mm"
$ws.Range("J24").Value = "1007.1 hPa"
$ws.Range("N24").Value = "11.0 °C 18:08 TU"
$ws.Range("E25").Value = "2026-02-11 18:49:30"
$ws.Range("I25").Value = "1.6 mm"
$ws.Range("E26").Value = "2026-02-11 18:49:33"
$ws.Range("H26").Value = "67%"
$ws.Range("J26").Value = "1003.0 hPa"
$ws.Range("L26").Value = "45.4 km/h - 168º 18:28 TU"
$ws.Range("E27").Value = "2026-02-11 18:49:35"
$ws.Range("I27").Value = "1.6 mm"
$ws.Range("E28").Value = "2026-02-11 18:49:38"
$ws.Range("H28").Value = "81%"
$ws.Range("J28").Value = "1003.3 hPa"
$ws.Range("O28").Value = "10.9 °C"
$ws.Range("E29").Value = "2026-02-11 18:49:41"
$ws.Range("E30").Value = "2026-02-11 18:49:43"
$ws.Range("J30").Value = "1003.2 hPa"
$ws.Range("E31").Value = "2026-02-11 18:49:46"
$ws.Range("J31").Value = "1002.4 hPa"
$ws.Range("E32").Value = "2026-02-11 18:49:49"
$ws.Range("H32").Value = "75%"
$ws.Range("I32").Value = "3.1 mm"
$ws.Range("E33").Value = "2026-02-11 18:49:51"
$ws.Range("J33").Value = "1004.9 hPa"
$ws.Range("E34").Value = "2026-02-11 18:49:54"
$ws.Range("O34").Value = "3.6 °C"
$ws.Range("E35").Value = "2026-02-11 18:49:57"
$ws.Range("J35").Value = "1007.8 hPa"
$ws.Range("O35").Value = "10.8 °C"
$ws.Range("E36").Value = "2026-02-11 18:49:59"
$ws.Range("J36").Value = "1003.4 hPa"
$ws.Range("E37").Value = "2026-02-11 18:50:02"
$ws.Range("J37").Value = "1004.5 hPa"
$ws.Range("E38").Value = "2026-02-11 18:50:05"
$ws.Range("H38").Value = "58%"
$ws.Range("E39").Value = "2026-02-11 18:50:07"
$ws.Range("H39").Value = "56%"
$ws.Range("E40").Value = "2026-02-11 18:50:10"
$ws.Range("H40").Value = "90%"
$ws.Range("I40").Value = "2.6 mm"
$ws.Range("J40").Value = "1007.1 hPa"
$ws.Range("E41").Value = "2026-02-11 18:50:13"
$ws.Range("H41").Value = "45%"
$ws.Range("J41").Value = "1004.9 hPa"
$ws.Range("N41").Value = "15.5 °C 18:29 TU"
$ws.Range("O41").Value = "19.2 °C"
$ws.Range("E42").Value = "2026-02-11 18:50:15"
$ws.Range("E43").Value = "2026-02-11 18:50:18"
$ws.Range("E44").Value = "2026-02-11 18:50:20"
$ws.Range("H44").Value = "83%"
$ws.Range("I44").Value = "5.8 mm"
$ws.Range("E45").Value = "2026-02-11 18:50:23"
$ws.Range("I45").Value = "2.8 mm"
$ws.Range("J45").Value = "1006.0 hPa"
$ws.Range("E46").Value = "2026-02-11 18:50:26"
$ws.Range("H46").Value = "58%"
$ws.Range("I46").Value = "1.6 mm"
$ws.Range("J46").Value = "1007.5 hPa"
$ws.Range("N46").Value = "12.8 °C 18:23 TU"
$ws.Range("O46").Value = "17.3 °C"
